# Replaces the fake "person" test data in rows 2-7 of the TestData sheet
# with a freshly (re)generated batch, matching the upstream fixture
# regeneration described by the commit ("fixed the excel issues").
#
# Columns: A=gender, B=first_name, C=Last_name, D=theDay, E=TheMonth,
#          F=theYear, G=email, H=company, I=testPassword
#
# theDay (column D) holds numeric-looking text (e.g. "30"); force text
# formatting before assigning so it round-trips as a string, not a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Vernon"
$ws.Range("C2").Value = "Rolfson"
$ws.Range("G2").Value = "yasmine.goodwin@hotmail.com"
$ws.Range("H2").Value = "Mayert, Cassin and Turcotte"
$ws.Range("I2").Value = "jq7bupm8ikj6"

# Row 3
$ws.Range("B3").Value = "Robert"
$ws.Range("C3").Value = "Brown"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30"
$ws.Range("G3").Value = "katelynn.bergnaum@gmail.com"
$ws.Range("H3").Value = "Becker LLC"
$ws.Range("I3").Value = "0e0b9lreyms3"

# Row 4
$ws.Range("B4").Value = "Cortez"
$ws.Range("C4").Value = "Olson"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "30"
$ws.Range("G4").Value = "casey.thompson@hotmail.com"
$ws.Range("H4").Value = "Moore Group"
$ws.Range("I4").Value = "7h62yjamdjf"

# Row 5
$ws.Range("B5").Value = "Cecil"
$ws.Range("C5").Value = "O'Kon"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "29"
$ws.Range("G5").Value = "shawnna.strosin@gmail.com"
$ws.Range("H5").Value = "Ruecker Group"
$ws.Range("I5").Value = "ue592kiios88uf"

# Row 6
$ws.Range("B6").Value = "Stacy"
$ws.Range("C6").Value = "Cassin"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "27"
$ws.Range("G6").Value = "dion.bogisich@gmail.com"
$ws.Range("H6").Value = "Thiel, Heathcote and Jerde"
$ws.Range("I6").Value = "im77t5d7rkimrf"

# Row 7
$ws.Range("B7").Value = "Lacy"
$ws.Range("C7").Value = "Krajcik"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "25"
$ws.Range("G7").Value = "newton.yost@hotmail.com"
$ws.Range("H7").Value = "Collins, Swaniawski and Dach"
$ws.Range("I7").Value = "bxy90qp7q3f7"
